# Apply updated crypto prices / 1h-volume percentages (cryptos list refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + "71.910.27"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.97%  "

$ws.Range("D3").Value = "'" + "4.045.56"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.96%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "'" + "540.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.76%  "

$ws.Range("D6").Value = "'" + "154.06"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +9.78%  "

$ws.Range("D7").Value = "'" + "0.696"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +14.55%  "

$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("D9").Value = "'" + "0.765"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +7.67%  "

$ws.Range("D10").Value = "'" + "0.174"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.62%  "

$ws.Range("E11").Value = "  +3.26%  "

$ws.Range("D12").Value = "'" + "48.40"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +16.66%  "

$ws.Range("E13").Value = "  +4.43%  "

$ws.Range("D14").Value = "'" + "4.691.83"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.00%  "

$ws.Range("D15").Value = "'" + "4.044.26"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.76%  "

$ws.Range("E16").Value = "  +2.53%  "

$ws.Range("D17").Value = "'" + "20.76"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.07%  "

$ws.Range("E18").Value = "  +2.29%  "

$ws.Range("E19").Value = "  -0.05%  "

$ws.Range("D20").Value = "'" + "71.909.96"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.96%  "

$ws.Range("D21").Value = "'" + "435.49"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.77%  "

$ws.Range("D22").Value = "'" + "99.41"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +14.53%  "

$ws.Range("D23").Value = "'" + "3.59"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.10%  "

$ws.Range("E24").Value = "  +7.48%  "

$ws.Range("D25").Value = "'" + "14.73"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.56%  "

$ws.Range("E26").Value = "  -3.07%  "

$ws.Range("D27").Value = "'" + "11.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.09%  "

$ws.Range("E28").Value = "  +4.79%  "

$ws.Range("D29").Value = "'" + "5.84"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.12%  "

$ws.Range("D30").Value = "'" + "3.64"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +30.47%  "

$ws.Range("D31").Value = "'" + "13.70"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.15%  "

$ws.Range("E32").Value = "  +5.93%  "

$ws.Range("D33").Value = "'" + "690.63"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.73%  "

$ws.Range("D34").Value = "'" + "6.97"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.72%  "

$ws.Range("D35").Value = "'" + "67.36"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.17%  "

$ws.Range("D36").Value = "'" + "43.31"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +9.99%  "

$ws.Range("E37").Value = "  -1.72%  "

$ws.Range("D38").Value = "'" + "0.157"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.32%  "

$ws.Range("E39").Value = "  +0.61%  "

$ws.Range("D40").Value = "'" + "3.47"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.31%  "

$ws.Range("D41").Value = "'" + "3.47"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +9.63%  "

$ws.Range("E42").Value = "  +0.00%  "

$ws.Range("E43").Value = "  +4.97%  "

$ws.Range("D44").Value = "'" + "0.998"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.20%  "

$ws.Range("E45").Value = "  +7.77%  "

$ws.Range("D46").Value = "'" + "2.73"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.03%  "

$ws.Range("D47").Value = "'" + "3.45"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.08%  "

$ws.Range("E48").Value = "  +10.00%  "

$ws.Range("E49").Value = "  +2.32%  "

$ws.Range("D50").Value = "'" + "3.36"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.91%  "

$ws.Range("E51").Value = "  -0.29%  "
